# The commit swaps the deck's two embedded themes: the "Integral" theme that
# currently drives the slide master (theme1.xml / Design #1) is replaced by
# the default "Office Theme" palette that used to sit unused in theme2.xml
# (wired only to the notes master). The font scheme (Arial-based "Office")
# and the fill/line/effect format scheme are already byte-identical between
# the two themes, so only the twelve theme colours actually need to change.

function Office-RGB($r, $g, $b) {
    # PowerPoint/VBA's RGB() packs colour components as 0xBBGGRR.
    return ($b * 65536) + ($g * 256) + $r
}

$p   = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index order exposed by ThemeColorScheme.Item(): dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - matches the <a:clrScheme> child order. Every
# slide shares the one slide-master theme, so touching slide 1 is enough.
$tcs.Item(1).RGB  = Office-RGB 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = Office-RGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = Office-RGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = Office-RGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = Office-RGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = Office-RGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = Office-RGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = Office-RGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = Office-RGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = Office-RGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = Office-RGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = Office-RGB 0x95 0x4F 0x72   # folHlink 954F72
